$wb = $excel.ActiveWorkbook

# 1. Insert a new worksheet "3D" before the first existing sheet (Architecture)
$firstSheet = $wb.Worksheets.Item(1)
$ws = $wb.Worksheets.Add($firstSheet)
$ws.Name = "3D"

# 2. Populate the "3D" sheet with data (header row + 11 article rows)
$ws.Cells.Item(1,1).Value = "更新日期"
$ws.Cells.Item(1,2).Value = "发表日期"
$ws.Cells.Item(1,3).Value = "文章标题"
$ws.Cells.Item(1,4).Value = "匹配关键词"
$ws.Cells.Item(1,5).Value = "作者"
$ws.Cells.Item(1,6).Value = "URL"
$ws.Cells.Item(1,7).Value = "摘要"
$ws.Cells.Item(2,1).Value = "2025-07-14"
$ws.Cells.Item(2,2).Value = "2024-12-11"
$ws.Cells.Item(2,3).Value = "SLGaussian: Fast Language Gaussian Splatting in Sparse Views"
$ws.Cells.Item(2,4).Value = "3DGS"
$ws.Cells.Item(2,5).Value = "Kangjie Chen, BingQuan Dai, Minghan Qin, Dongbin Zhang, Peihao Li, Yingshuang Zou, Haoqian Wang"
$ws.Cells.Item(2,6).Value = "http://arxiv.org/abs/2412.08331v2"
$ws.Cells.Item(2,7).Value = "3D semantic field learning is crucial for applications like autonomous
navigation, AR/VR, and robotics, where accurate comprehension of 3D scenes from
limited viewpoints is essential. Existing methods struggle under sparse view
conditions, relying on inefficient per-scene multi-view optimizations, which
are impractical for many real-world tasks. To address this, we propose
SLGaussian, a feed-forward method for constructing 3D semantic fields from
sparse viewpoints, allowing direct inference of 3DGS-based scenes. By ensuring
consistent SAM segmentations through video tracking and using low-dimensional
indexing for high-dimensional CLIP features, SLGaussian efficiently embeds
language information in 3D space, offering a robust solution for accurate 3D
scene understanding under sparse view conditions. In experiments on two-view
sparse 3D object querying and segmentation in the LERF and 3D-OVS datasets,
SLGaussian outperforms existing methods in chosen IoU, Localization Accuracy,
and mIoU. Moreover, our model achieves scene inference in under 30 seconds and
open-vocabulary querying in just 0.011 seconds per query."
$ws.Cells.Item(3,1).Value = "2025-07-14"
$ws.Cells.Item(3,2).Value = "2025-03-07"
$ws.Cells.Item(3,3).Value = "CoMoGaussian: Continuous Motion-Aware Gaussian Splatting from Motion-Blurred Images"
$ws.Cells.Item(3,4).Value = "3DGS"
$ws.Cells.Item(3,5).Value = "Jungho Lee, Donghyeong Kim, Dogyoon Lee, Suhwan Cho, Minhyeok Lee, Wonjoon Lee, Taeoh Kim, Dongyoon Wee, Sangyoun Lee"
$ws.Cells.Item(3,6).Value = "http://arxiv.org/abs/2503.05332v2"
$ws.Cells.Item(3,7).Value = "3D Gaussian Splatting (3DGS) has gained significant attention due to its
high-quality novel view rendering, motivating research to address real-world
challenges. A critical issue is the camera motion blur caused by movement
during exposure, which hinders accurate 3D scene reconstruction. In this study,
we propose CoMoGaussian, a Continuous Motion-Aware Gaussian Splatting that
reconstructs precise 3D scenes from motion-blurred images while maintaining
real-time rendering speed. Considering the complex motion patterns inherent in
real-world camera movements, we predict continuous camera trajectories using
neural ordinary differential equations (ODEs). To ensure accurate modeling, we
employ rigid body transformations, preserving the shape and size of the object
but rely on the discrete integration of sampled frames. To better approximate
the continuous nature of motion blur, we introduce a continuous motion
refinement (CMR) transformation that refines rigid transformations by
incorporating additional learnable parameters. By revisiting fundamental camera
theory and leveraging advanced neural ODE techniques, we achieve precise
modeling of continuous camera trajectories, leading to improved reconstruction
accuracy. Extensive experiments demonstrate state-of-the-art performance both
quantitatively and qualitatively on benchmark datasets, which include a wide
range of motion blur scenarios, from moderate to extreme blur."
$ws.Cells.Item(4,1).Value = "2025-07-14"
$ws.Cells.Item(4,2).Value = "2024-09-17"
$ws.Cells.Item(4,3).Value = "HGSLoc: 3DGS-based Heuristic Camera Pose Refinement"
$ws.Cells.Item(4,4).Value = "NeRF, Neural Rendering"
$ws.Cells.Item(4,5).Value = "Zhongyan Niu, Zhen Tan, Jinpu Zhang, Xueliang Yang, Dewen Hu"
$ws.Cells.Item(4,6).Value = "http://arxiv.org/abs/2409.10925v3"
$ws.Cells.Item(4,7).Value = "Visual localization refers to the process of determining camera poses and
orientation within a known scene representation. This task is often complicated
by factors such as changes in illumination and variations in viewing angles. In
this paper, we propose HGSLoc, a novel lightweight plug-and-play pose
optimization framework, which integrates 3D reconstruction with a heuristic
refinement strategy to achieve higher pose estimation accuracy. Specifically,
we introduce an explicit geometric map for 3D representation and high-fidelity
rendering, allowing the generation of high-quality synthesized views to support
accurate visual localization. Our method demonstrates higher localization
accuracy compared to NeRF-based neural rendering localization approaches. We
introduce a heuristic refinement strategy, its efficient optimization
capability can quickly locate the target node, while we set the step level
optimization step to enhance the pose accuracy in the scenarios with small
errors. With carefully designed heuristic functions, it offers efficient
optimization capabilities, enabling rapid error reduction in rough localization
estimations. Our method mitigates the dependence on complex neural network
models while demonstrating improved robustness against noise and higher
localization accuracy in challenging environments, as compared to neural
network joint optimization strategies. The optimization framework proposed in
this paper introduces novel approaches to visual localization by integrating
the advantages of 3D reconstruction and the heuristic refinement strategy,
which demonstrates strong performance across multiple benchmark datasets,
including 7Scenes and Deep Blending dataset. The implementation of our method
has been released at https://github.com/anchang699/HGSLoc."
$ws.Cells.Item(5,1).Value = "2025-07-14"
$ws.Cells.Item(5,2).Value = "2025-07-14"
$ws.Cells.Item(5,3).Value = "3DGAA: Realistic and Robust 3D Gaussian-based Adversarial Attack for Autonomous Driving"
$ws.Cells.Item(5,4).Value = "3DGS"
$ws.Cells.Item(5,5).Value = "Yixun Zhang, Lizhi Wang, Junjun Zhao, Wending Zhao, Feng Zhou, Yonghao Dang, Jianqin Yin"
$ws.Cells.Item(5,6).Value = "http://arxiv.org/abs/2507.09993v1"
$ws.Cells.Item(5,7).Value = "Camera-based object detection systems play a vital role in autonomous
driving, yet they remain vulnerable to adversarial threats in real-world
environments. While existing 2D and 3D physical attacks typically optimize
texture, they often struggle to balance physical realism and attack robustness.
In this work, we propose 3D Gaussian-based Adversarial Attack (3DGAA), a novel
adversarial object generation framework that leverages the full 14-dimensional
parameterization of 3D Gaussian Splatting (3DGS) to jointly optimize geometry
and appearance in physically realizable ways. Unlike prior works that rely on
patches or texture, 3DGAA jointly perturbs both geometric attributes (shape,
scale, rotation) and appearance attributes (color, opacity) to produce
physically realistic and transferable adversarial objects. We further introduce
a physical filtering module to preserve geometric fidelity, and a physical
augmentation module to simulate complex physical scenarios, thus enhancing
attack generalization under real-world conditions. We evaluate 3DGAA on both
virtual benchmarks and physical-world setups using miniature vehicle models.
Experimental results show that 3DGAA achieves to reduce the detection mAP from
87.21% to 7.38%, significantly outperforming existing 3D physical attacks.
Moreover, our method maintains high transferability across different physical
conditions, demonstrating a new state-of-the-art in physically realizable
adversarial attacks. These results validate 3DGAA as a practical attack
framework for evaluating the safety of perception systems in autonomous
driving."
$ws.Cells.Item(6,1).Value = "2025-07-14"
$ws.Cells.Item(6,2).Value = "2025-07-11"
$ws.Cells.Item(6,3).Value = "CLiFT: Compressive Light-Field Tokens for Compute-Efficient and Adaptive Neural Rendering"
$ws.Cells.Item(6,4).Value = "Neural Rendering"
$ws.Cells.Item(6,5).Value = "Zhengqing Wang, Yuefan Wu, Jiacheng Chen, Fuyang Zhang, Yasutaka Furukawa"
$ws.Cells.Item(6,6).Value = "http://arxiv.org/abs/2507.08776v2"
$ws.Cells.Item(6,7).Value = "This paper proposes a neural rendering approach that represents a scene as
""compressed light-field tokens (CLiFTs)"", retaining rich appearance and
geometric information of a scene. CLiFT enables compute-efficient rendering by
compressed tokens, while being capable of changing the number of tokens to
represent a scene or render a novel view with one trained network. Concretely,
given a set of images, multi-view encoder tokenizes the images with the camera
poses. Latent-space K-means selects a reduced set of rays as cluster centroids
using the tokens. The multi-view ````condenser'' compresses the information of
all the tokens into the centroid tokens to construct CLiFTs. At test time,
given a target view and a compute budget (i.e., the number of CLiFTs), the
system collects the specified number of nearby tokens and synthesizes a novel
view using a compute-adaptive renderer. Extensive experiments on RealEstate10K
and DL3DV datasets quantitatively and qualitatively validate our approach,
achieving significant data reduction with comparable rendering quality and the
highest overall rendering score, while providing trade-offs of data size,
rendering quality, and rendering speed."
$ws.Cells.Item(7,1).Value = "2025-07-12"
$ws.Cells.Item(7,2).Value = "2025-07-12"
$ws.Cells.Item(7,3).Value = "Stable Score Distillation"
$ws.Cells.Item(7,4).Value = "NeRF"
$ws.Cells.Item(7,5).Value = "Haiming Zhu, Yangyang Xu, Chenshu Xu, Tingrui Shen, Wenxi Liu, Yong Du, Jun Yu, Shengfeng He"
$ws.Cells.Item(7,6).Value = "http://arxiv.org/abs/2507.09168v1"
$ws.Cells.Item(7,7).Value = "Text-guided image and 3D editing have advanced with diffusion-based models,
yet methods like Delta Denoising Score often struggle with stability, spatial
control, and editing strength. These limitations stem from reliance on complex
auxiliary structures, which introduce conflicting optimization signals and
restrict precise, localized edits. We introduce Stable Score Distillation
(SSD), a streamlined framework that enhances stability and alignment in the
editing process by anchoring a single classifier to the source prompt.
Specifically, SSD utilizes Classifier-Free Guidance (CFG) equation to achieves
cross-prompt alignment, and introduces a constant term null-text branch to
stabilize the optimization process. This approach preserves the original
content's structure and ensures that editing trajectories are closely aligned
with the source prompt, enabling smooth, prompt-specific modifications while
maintaining coherence in surrounding regions. Additionally, SSD incorporates a
prompt enhancement branch to boost editing strength, particularly for style
transformations. Our method achieves state-of-the-art results in 2D and 3D
editing tasks, including NeRF and text-driven style edits, with faster
convergence and reduced complexity, providing a robust and efficient solution
for text-guided editing."
$ws.Cells.Item(8,1).Value = "2025-07-11"
$ws.Cells.Item(8,2).Value = "2025-07-11"
$ws.Cells.Item(8,3).Value = "From images to properties: a NeRF-driven framework for granular material parameter inversion"
$ws.Cells.Item(8,4).Value = "NeRF"
$ws.Cells.Item(8,5).Value = "Cheng-Hsi Hsiao, Krishna Kumar"
$ws.Cells.Item(8,6).Value = "http://arxiv.org/abs/2507.09005v1"
$ws.Cells.Item(8,7).Value = "We introduce a novel framework that integrates Neural Radiance Fields (NeRF)
with Material Point Method (MPM) simulation to infer granular material
properties from visual observations. Our approach begins by generating
synthetic experimental data, simulating an plow interacting with sand. The
experiment is rendered into realistic images as the photographic observations.
These observations include multi-view images of the experiment's initial state
and time-sequenced images from two fixed cameras. Using NeRF, we reconstruct
the 3D geometry from the initial multi-view images, leveraging its capability
to synthesize novel viewpoints and capture intricate surface details. The
reconstructed geometry is then used to initialize material point positions for
the MPM simulation, where the friction angle remains unknown. We render images
of the simulation under the same camera setup and compare them to the observed
images. By employing Bayesian optimization, we minimize the image loss to
estimate the best-fitting friction angle. Our results demonstrate that friction
angle can be estimated with an error within 2 degrees, highlighting the
effectiveness of inverse analysis through purely visual observations. This
approach offers a promising solution for characterizing granular materials in
real-world scenarios where direct measurement is impractical or impossible."
$ws.Cells.Item(9,1).Value = "2025-07-11"
$ws.Cells.Item(9,2).Value = "2024-10-31"
$ws.Cells.Item(9,3).Value = "GeoSplatting: Towards Geometry Guided Gaussian Splatting for Physically-based Inverse Rendering"
$ws.Cells.Item(9,4).Value = "3DGS"
$ws.Cells.Item(9,5).Value = "Kai Ye, Chong Gao, Guanbin Li, Wenzheng Chen, Baoquan Chen"
$ws.Cells.Item(9,6).Value = "http://arxiv.org/abs/2410.24204v3"
$ws.Cells.Item(9,7).Value = "Recent 3D Gaussian Splatting (3DGS) representations have demonstrated
remarkable performance in novel view synthesis; further, material-lighting
disentanglement on 3DGS warrants relighting capabilities and its adaptability
to broader applications. While the general approach to the latter operation
lies in integrating differentiable physically-based rendering (PBR) techniques
to jointly recover BRDF materials and environment lighting, achieving a precise
disentanglement remains an inherently difficult task due to the challenge of
accurately modeling light transport. Existing approaches typically approximate
Gaussian points' normals, which constitute an implicit geometric constraint.
However, they usually suffer from inaccuracies in normal estimation that
subsequently degrade light transport, resulting in noisy material decomposition
and flawed relighting results. To address this, we propose GeoSplatting, a
novel approach that augments 3DGS with explicit geometry guidance for precise
light transport modeling. By differentiably constructing a surface-grounded
3DGS from an optimizable mesh, our approach leverages well-defined mesh normals
and the opaque mesh surface, and additionally facilitates the use of mesh-based
ray tracing techniques for efficient, occlusion-aware light transport
calculations. This enhancement ensures precise material decomposition while
preserving the efficiency and high-quality rendering capabilities of 3DGS.
Comprehensive evaluations across diverse datasets demonstrate the effectiveness
of GeoSplatting, highlighting its superior efficiency and state-of-the-art
inverse rendering performance. The project page can be found at
https://pku-vcl-geometry.github.io/GeoSplatting/."
$ws.Cells.Item(10,1).Value = "2025-07-10"
$ws.Cells.Item(10,2).Value = "2025-07-10"
$ws.Cells.Item(10,3).Value = "RegGS: Unposed Sparse Views Gaussian Splatting with 3DGS Registration"
$ws.Cells.Item(10,4).Value = "3DGS"
$ws.Cells.Item(10,5).Value = "Chong Cheng, Yu Hu, Sicheng Yu, Beizhen Zhao, Zijian Wang, Hao Wang"
$ws.Cells.Item(10,6).Value = "http://arxiv.org/abs/2507.08136v1"
$ws.Cells.Item(10,7).Value = "3D Gaussian Splatting (3DGS) has demonstrated its potential in reconstructing
scenes from unposed images. However, optimization-based 3DGS methods struggle
with sparse views due to limited prior knowledge. Meanwhile, feed-forward
Gaussian approaches are constrained by input formats, making it challenging to
incorporate more input views. To address these challenges, we propose RegGS, a
3D Gaussian registration-based framework for reconstructing unposed sparse
views. RegGS aligns local 3D Gaussians generated by a feed-forward network into
a globally consistent 3D Gaussian representation. Technically, we implement an
entropy-regularized Sinkhorn algorithm to efficiently solve the optimal
transport Mixture 2-Wasserstein `$(\text{MW}_2)`$ distance, which serves as an
alignment metric for Gaussian mixture models (GMMs) in `$\mathrm{Sim}(3)`$ space.
Furthermore, we design a joint 3DGS registration module that integrates the
`$\text{MW}_2`$ distance, photometric consistency, and depth geometry. This
enables a coarse-to-fine registration process while accurately estimating
camera poses and aligning the scene. Experiments on the RE10K and ACID datasets
demonstrate that RegGS effectively registers local Gaussians with high
fidelity, achieving precise pose estimation and high-quality novel-view
synthesis. Project page: https://3dagentworld.github.io/reggs/."
$ws.Cells.Item(11,1).Value = "2025-07-10"
$ws.Cells.Item(11,2).Value = "2025-07-10"
$ws.Cells.Item(11,3).Value = "RTR-GS: 3D Gaussian Splatting for Inverse Rendering with Radiance Transfer and Reflection"
$ws.Cells.Item(11,4).Value = "3DGS"
$ws.Cells.Item(11,5).Value = "Yongyang Zhou, Fang-Lue Zhang, Zichen Wang, Lei Zhang"
$ws.Cells.Item(11,6).Value = "http://arxiv.org/abs/2507.07733v1"
$ws.Cells.Item(11,7).Value = "3D Gaussian Splatting (3DGS) has demonstrated impressive capabilities in
novel view synthesis. However, rendering reflective objects remains a
significant challenge, particularly in inverse rendering and relighting. We
introduce RTR-GS, a novel inverse rendering framework capable of robustly
rendering objects with arbitrary reflectance properties, decomposing BRDF and
lighting, and delivering credible relighting results. Given a collection of
multi-view images, our method effectively recovers geometric structure through
a hybrid rendering model that combines forward rendering for radiance transfer
with deferred rendering for reflections. This approach successfully separates
high-frequency and low-frequency appearances, mitigating floating artifacts
caused by spherical harmonic overfitting when handling high-frequency details.
We further refine BRDF and lighting decomposition using an additional
physically-based deferred rendering branch. Experimental results show that our
method enhances novel view synthesis, normal estimation, decomposition, and
relighting while maintaining efficient training inference process."
$ws.Cells.Item(12,1).Value = "2025-07-10"
$ws.Cells.Item(12,2).Value = "2025-07-10"
$ws.Cells.Item(12,3).Value = "MUVOD: A Novel Multi-view Video Object Segmentation Dataset and A Benchmark for 3D Segmentation"
$ws.Cells.Item(12,4).Value = "NeRF"
$ws.Cells.Item(12,5).Value = "Bangning Wei, Joshua Maraval, Meriem Outtas, Kidiyo Kpalma, Nicolas Ramin, Lu Zhang"
$ws.Cells.Item(12,6).Value = "http://arxiv.org/abs/2507.07519v1"
$ws.Cells.Item(12,7).Value = "The application of methods based on Neural Radiance Fields (NeRF) and 3D
Gaussian Splatting (3D GS) have steadily gained popularity in the field of 3D
object segmentation in static scenes. These approaches demonstrate efficacy in
a range of 3D scene understanding and editing tasks. Nevertheless, the 4D
object segmentation of dynamic scenes remains an underexplored field due to the
absence of a sufficiently extensive and accurately labelled multi-view video
dataset. In this paper, we present MUVOD, a new multi-view video dataset for
training and evaluating object segmentation in reconstructed real-world
scenarios. The 17 selected scenes, describing various indoor or outdoor
activities, are collected from different sources of datasets originating from
various types of camera rigs. Each scene contains a minimum of 9 views and a
maximum of 46 views. We provide 7830 RGB images (30 frames per video) with
their corresponding segmentation mask in 4D motion, meaning that any object of
interest in the scene could be tracked across temporal frames of a given view
or across different views belonging to the same camera rig. This dataset, which
contains 459 instances of 73 categories, is intended as a basic benchmark for
the evaluation of multi-view video segmentation methods. We also present an
evaluation metric and a baseline segmentation approach to encourage and
evaluate progress in this evolving field. Additionally, we propose a new
benchmark for 3D object segmentation task with a subset of annotated multi-view
images selected from our MUVOD dataset. This subset contains 50 objects of
different conditions in different scenarios, providing a more comprehensive
analysis of state-of-the-art 3D object segmentation methods. Our proposed MUVOD
dataset is available at https://volumetric-repository.labs.b-com.com/#/muvod."

# 3. Style the header row (bold, centered, top-aligned, thin border) to match
#    the look of the header rows on the other sheets.
$header = $ws.Range("A1:G1")
$header.Font.Bold = $true
$header.HorizontalAlignment = -4108
$header.VerticalAlignment = -4160
$header.Borders.LineStyle = 1

# 4. Rename "LLM Other Compression" to "LLM Compression Others"
$wsCompression = $wb.Worksheets.Item("LLM Other Compression")
$wsCompression.Name = "LLM Compression Others"
